$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.956140350877193
$ws.Range("C2").Value = 0.8859649122807017
$ws.Range("D2").Value = 0.9385964912280702
$ws.Range("E2").Value = 0.9298245614035088
$ws.Range("F2").Value = 0.956140350877193
$ws.Range("G2").Value = 0.9298245614035088
$ws.Range("H2").Value = 0.956140350877193
